# Update vm_pu results for Case_1_23 (380 kV case)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.033048350214858
$ws.Range("D2").Value = 1.049936085579724
$ws.Range("E2").Value = 1.043392414093089
$ws.Range("F2").Value = 1.056693077107794
$ws.Range("I2").Value = 1.041599793834862
$ws.Range("J2").Value = 1.03817521648763
$ws.Range("K2").Value = 1.052691758345476
$ws.Range("L2").Value = 1.046166408771222
$ws.Range("M2").Value = 1.059430088694119
$ws.Range("N2").Value = 1.016701521056828

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.033809436172126
$ws.Range("D3").Value = 1.050478445362427
$ws.Range("E3").Value = 1.044044506942185
$ws.Range("F3").Value = 1.0573685816333
$ws.Range("I3").Value = 1.041744946503135
$ws.Range("J3").Value = 1.038579809084225
$ws.Range("K3").Value = 1.053047160829539
$ws.Range("L3").Value = 1.046629943025823
$ws.Range("M3").Value = 1.05991963423988
$ws.Range("N3").Value = 1.016836773943497

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.034302767882629
$ws.Range("D4").Value = 1.050830023517636
$ws.Range("E4").Value = 1.04446756476127
$ws.Range("F4").Value = 1.057806739653426
$ws.Range("I4").Value = 1.041838097727256
$ws.Range("J4").Value = 1.038841770502437
$ws.Range("K4").Value = 1.05327705418585
$ws.Range("L4").Value = 1.046930302783533
$ws.Range("M4").Value = 1.060236780014811
$ws.Range("N4").Value = 1.016924315018851

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.034510367969049
$ws.Range("D5").Value = 1.050977976941212
$ws.Range("E5").Value = 1.044645681917386
$ws.Range("F5").Value = 1.057991193119305
$ws.Range("I5").Value = 1.041877072809386
$ws.Range("J5").Value = 1.038951936708822
$ws.Range("K5").Value = 1.053373681719965
$ws.Range("L5").Value = 1.047056673479358
$ws.Range("M5").Value = 1.060370196593117
$ws.Range("N5").Value = 1.016961122392146

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.034545236807876
$ws.Range("D6").Value = 1.051002827705643
$ws.Range("E6").Value = 1.044675603970403
$ws.Range("F6").Value = 1.058022178369385
$ws.Range("I6").Value = 1.041883605984481
$ws.Range("J6").Value = 1.038970436265748
$ws.Range("K6").Value = 1.0533899047137
$ws.Range("L6").Value = 1.04707789746076
$ws.Range("M6").Value = 1.060392602961312
$ws.Range("N6").Value = 1.016967302796903

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.034305541048315
$ws.Range("D7").Value = 1.050831999891878
$ws.Range("E7").Value = 1.044469943738141
$ws.Range("F7").Value = 1.057809203342993
$ws.Range("I7").Value = 1.041838619244869
$ws.Range("J7").Value = 1.038843242402574
$ws.Range("K7").Value = 1.053278345406699
$ws.Range("L7").Value = 1.046931990966182
$ws.Range("M7").Value = 1.060238562388175
$ws.Range("N7").Value = 1.016924806821343

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.033305383865525
$ws.Range("D8").Value = 1.050119245855258
$ws.Range("E8").Value = 1.043612560788099
$ws.Range("F8").Value = 1.05692114588517
$ws.Range("I8").Value = 1.041649008271236
$ws.Range("J8").Value = 1.038311915876606
$ws.Range("K8").Value = 1.052811883086454
$ws.Range("L8").Value = 1.046322973984301
$ws.Range("M8").Value = 1.059595453760283
$ws.Range("N8").Value = 1.016747225254732

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.031549641428315
$ws.Range("D9").Value = 1.048868238231785
$ws.Range("E9").Value = 1.042110338619025
$ws.Range("F9").Value = 1.055364504875898
$ws.Range("I9").Value = 1.041309010926926
$ws.Range("J9").Value = 1.037376959972023
$ws.Range("K9").Value = 1.051989401769398
$ws.Range("L9").Value = 1.045253111683132
$ws.Range("M9").Value = 1.058465180698677
$ws.Range("N9").Value = 1.016434505274183

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.030383746951402
$ws.Range("D10").Value = 1.048037689823383
$ws.Range("E10").Value = 1.041114763783994
$ws.Range("F10").Value = 1.054332413528241
$ws.Range("I10").Value = 1.041078442389534
$ws.Range("J10").Value = 1.036754619339119
$ws.Range("K10").Value = 1.051440823377287
$ws.Range("L10").Value = 1.044542186896776
$ws.Range("M10").Value = 1.05771376631326
$ws.Range("N10").Value = 1.01622619064853

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.02988001864528
$ws.Range("D11").Value = 1.04767889954844
$ws.Range("E11").Value = 1.040685096876186
$ws.Range("F11").Value = 1.053886880223416
$ws.Range("I11").Value = 1.040977688135388
$ws.Range("J11").Value = 1.036485385078916
$ws.Range("K11").Value = 1.051203240674697
$ws.Range("L11").Value = 1.044234918453614
$ws.Range("M11").Value = 1.057388916547737
$ws.Range("N11").Value = 1.016136033691413

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.029693080636446
$ws.Range("D12").Value = 1.04754575763448
$ws.Range("E12").Value = 1.040525715545588
$ws.Range("F12").Value = 1.053721597375656
$ws.Range("I12").Value = 1.040940126649484
$ws.Range("J12").Value = 1.036385417480344
$ws.Range("K12").Value = 1.051114986711069
$ws.Range("L12").Value = 1.044120872070487
$ws.Range("M12").Value = 1.05726833249583
$ws.Range("N12").Value = 1.016102552608464

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.029733171787604
$ws.Range("D13").Value = 1.047574311165155
$ws.Range("E13").Value = 1.040559893567994
$ws.Range("F13").Value = 1.053757041654404
$ws.Range("I13").Value = 1.040948189895993
$ws.Range("J13").Value = 1.036406859127296
$ws.Range("K13").Value = 1.051133917695136
$ws.Range("L13").Value = 1.044145331452394
$ws.Range("M13").Value = 1.057294194551914
$ws.Range("N13").Value = 1.016109734079166

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.029864562829201
$ws.Range("D14").Value = 1.047667891356928
$ws.Range("E14").Value = 1.040671917944718
$ws.Range("F14").Value = 1.053873213633612
$ws.Range("I14").Value = 1.040974586078532
$ws.Range("J14").Value = 1.036477120942692
$ws.Range("K14").Value = 1.051195945671111
$ws.Range("L14").Value = 1.044225489563767
$ws.Range("M14").Value = 1.057378947397594
$ws.Range("N14").Value = 1.016133265984274

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.029945539658057
$ws.Range("D15").Value = 1.047725566329542
$ws.Range("E15").Value = 1.040740968560674
$ws.Range("F15").Value = 1.053944818664901
$ws.Range("I15").Value = 1.040990831526109
$ws.Range("J15").Value = 1.036520416643447
$ws.Range("K15").Value = 1.051234162515205
$ws.Range("L15").Value = 1.044274889172904
$ws.Range("M15").Value = 1.057431177027988
$ws.Range("N15").Value = 1.016147765740514

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.030417201340861
$ws.Range("D16").Value = 1.048061519491511
$ws.Range("E16").Value = 1.041143309543683
$ws.Range("F16").Value = 1.054362011153564
$ws.Range("I16").Value = 1.041085109877804
$ws.Range("J16").Value = 1.036772492779623
$ws.Range("K16").Value = 1.051456590149566
$ws.Range("L16").Value = 1.04456259137249
$ws.Range("M16").Value = 1.057735336585857
$ws.Range("N16").Value = 1.016232175053728

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.030713361362913
$ws.Range("D17").Value = 1.048272481172466
$ws.Range("E17").Value = 1.041396070106799
$ws.Range("F17").Value = 1.054624073192667
$ws.Range("I17").Value = 1.041144003396271
$ws.Range("J17").Value = 1.036930679642094
$ws.Range("K17").Value = 1.051596102262881
$ws.Range("L17").Value = 1.044743212250778
$ws.Range("M17").Value = 1.057926267728437
$ws.Range("N17").Value = 1.016285135105215

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.030886213615281
$ws.Range("D18").Value = 1.048395612711806
$ws.Range("E18").Value = 1.041543638272015
$ws.Range("F18").Value = 1.054777061481827
$ws.Range("I18").Value = 1.041178266498971
$ws.Range("J18").Value = 1.037022970745402
$ws.Range("K18").Value = 1.051677472908295
$ws.Range("L18").Value = 1.044848619918478
$ws.Range("M18").Value = 1.058037684452847
$ws.Range("N18").Value = 1.016316030080733

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.030945169872558
$ws.Range("D19").Value = 1.048437611084818
$ws.Range("E19").Value = 1.041593978407756
$ws.Range("F19").Value = 1.054829248841131
$ws.Range("I19").Value = 1.041189934296767
$ws.Range("J19").Value = 1.037054443549082
$ws.Range("K19").Value = 1.051705217400778
$ws.Range("L19").Value = 1.044884570413847
$ws.Range("M19").Value = 1.05807568307054
$ws.Range("N19").Value = 1.016326565178108

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.030681575100153
$ws.Range("D20").Value = 1.048249838570227
$ws.Range("E20").Value = 1.041368937082741
$ws.Range("F20").Value = 1.054595942752188
$ws.Range("I20").Value = 1.041137693826067
$ws.Range("J20").Value = 1.036913705259894
$ws.Range("K20").Value = 1.051581134378514
$ws.Range("L20").Value = 1.044723827684007
$ws.Range("M20").Value = 1.057905777457434
$ws.Range("N20").Value = 1.016279452547145

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.02982586675714
$ws.Range("D21").Value = 1.047640330750831
$ws.Range("E21").Value = 1.040638923575002
$ws.Range("F21").Value = 1.053838998118764
$ws.Range("I21").Value = 1.040966816832271
$ws.Range("J21").Value = 1.036456429528975
$ws.Range("K21").Value = 1.051177680108568
$ws.Range("L21").Value = 1.044201882588491
$ws.Range("M21").Value = 1.057353987584363
$ws.Range("N21").Value = 1.01612633622114

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.029288828428862
$ws.Range("D22").Value = 1.047257855081689
$ws.Range("E22").Value = 1.040181186431007
$ws.Range("F22").Value = 1.053364281921081
$ws.Range("I22").Value = 1.04085858821011
$ws.Range("J22").Value = 1.036169142654264
$ws.Range("K22").Value = 1.050923983245567
$ws.Range("L22").Value = 1.043874217996169
$ws.Range("M22").Value = 1.05700751693059
$ws.Range("N22").Value = 1.016030107967643

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.029573429066209
$ws.Range("D23").Value = 1.047460541227653
$ws.Range("E23").Value = 1.040423722202343
$ws.Range("F23").Value = 1.053615822954206
$ws.Range("I23").Value = 1.04091603699493
$ws.Range("J23").Value = 1.036321417480305
$ws.Range("K23").Value = 1.051058475003374
$ws.Range("L23").Value = 1.044047870963371
$ws.Range("M23").Value = 1.057191143214703
$ws.Range("N23").Value = 1.016081116226532

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.030695937613366
$ws.Range("D24").Value = 1.04826006953707
$ws.Range("E24").Value = 1.041381196904753
$ws.Range("F24").Value = 1.054608653280052
$ws.Range("I24").Value = 1.041140545122785
$ws.Range("J24").Value = 1.036921375180505
$ws.Range("K24").Value = 1.051587897735012
$ws.Range("L24").Value = 1.044732586566963
$ws.Range("M24").Value = 1.05791503597523
$ws.Range("N24").Value = 1.0162820202374

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.032002740940922
$ws.Range("D25").Value = 1.049191053522649
$ws.Range("E25").Value = 1.042497667264854
$ws.Range("F25").Value = 1.055765944583851
$ws.Range("I25").Value = 1.041397600009339
$ws.Range("J25").Value = 1.037618504767872
$ws.Range("K25").Value = 1.052202084787966
$ws.Range("L25").Value = 1.045529295388852
$ws.Range("M25").Value = 1.058757020849704
$ws.Range("N25").Value = 1.016515323927864
